# Generate Report for Handback
#
# For both the "zh-cn" and "de-de" language sheets, row 8 (the
# e8a43dbf-4a30-4d44-a948-ac925269ecd4 entry) gets a freshly generated
# handback-status report:
#   - a new hyperlinked "Latest Target File" (column I) pointing at the
#     latest handback commit
#   - an updated "Latest Handback DateTime" (column K)
#   - a new "Error Detail" (column P) explaining the handback file is
#     stale, with the current vs. latest commit URLs
#   - the de-de sheet additionally refreshes "Latest Handback File" (J)
# The "Error Detail" column (P, the 16th column) is widened on both
# sheets to fit the new message.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6079a4471f2bd8a8e9df28aadf363e406c74f500/e2e/e8a43dbf-4a30-4d44-a948-ac925269ecd4.md"
$displayName = "e8a43dbf-4a30-4d44-a948-ac925269ecd4.md"
$errorDetail = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/021a50e29e9389c83fb036076e4bd1f414437e5f/e2e/e8a43dbf-4a30-4d44-a948-ac925269ecd4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6079a4471f2bd8a8e9df28aadf363e406c74f500/e2e/e8a43dbf-4a30-4d44-a948-ac925269ecd4.md.'

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("P1").ColumnWidth = 39.17

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $latestUrl, "", "", $displayName)
$wsZh.Range("J8").Value = "e8a43dbf-4a30-4d44-a948-ac925269ecd4.c7a6959c0741eb1a0ddc9862c114ced5ee89af3b.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-25 02:42:32"
$wsZh.Range("P8").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("P1").ColumnWidth = 39.17

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $latestUrl, "", "", $displayName)
$wsDe.Range("J8").Value = "e8a43dbf-4a30-4d44-a948-ac925269ecd4.c7a6959c0741eb1a0ddc9862c114ced5ee89af3b.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-25 02:42:40"
$wsDe.Range("P8").Value = $errorDetail
